$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3199.5
$ws.Range("I2").Value = 2047.5883
$ws.Range("J2").Value = 7116
$ws.Range("K2").Value = 2047.5883
$ws.Range("L2").Value = 7116
$ws.Range("M2").Value = -1934.5883
$ws.Range("N2").Value = -7342

$ws.Range("H32").Value = 2827864.2
$ws.Range("I32").Value = 4912.26
$ws.Range("J32").Value = 20941806
$ws.Range("K32").Value = 4912.26
$ws.Range("L32").Value = 20941806
$ws.Range("M32").Value = -4625.26
$ws.Range("N32").Value = -20942380

$ws.Range("H45").Value = 3182.6206
$ws.Range("I45").Value = 2375.9375
$ws.Range("J45").Value = 4175.4614
$ws.Range("K45").Value = 2375.9375
$ws.Range("L45").Value = 4175.4614
$ws.Range("M45").Value = -1998.9375
$ws.Range("N45").Value = -4929.4614

$ws.Range("H61").Value = 1406.921
$ws.Range("I61").Value = 1390.8918
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1390.8918
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -1178.8918
$ws.Range("N61").Value = -2424

$ws.Range("H74").Value = 827.48334
$ws.Range("I74").Value = 818.35297
$ws.Range("J74").Value = 879.2222
$ws.Range("K74").Value = 818.35297
$ws.Range("L74").Value = 879.2222
$ws.Range("M74").Value = 55.64702999999997
$ws.Range("N74").Value = -2627.2222

$ws.Range("H77").Value = 827.48334
$ws.Range("I77").Value = 818.35297
$ws.Range("J77").Value = 879.2222
$ws.Range("K77").Value = 4091.76485
$ws.Range("L77").Value = 4396.111
$ws.Range("M77").Value = 276.23515
$ws.Range("N77").Value = -13132.111

$ws.Range("H116").Value = 3199.5
$ws.Range("I116").Value = 2047.5883
$ws.Range("J116").Value = 7116
$ws.Range("K116").Value = 2047.5883
$ws.Range("L116").Value = 7116
$ws.Range("M116").Value = 246.4117000000001
$ws.Range("N116").Value = -11704

$ws.Range("H136").Value = 1406.921
$ws.Range("I136").Value = 1390.8918
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 4172.6754
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -1622.6754
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3199.5
$ws.Range("I3").Value = 2047.5883
$ws.Range("J3").Value = 7116
$ws.Range("K3").Value = 2047.5883
$ws.Range("L3").Value = 7116
$ws.Range("M3").Value = -1933.5883
$ws.Range("N3").Value = -7344

$ws.Range("H80").Value = 412.0435
$ws.Range("J80").Value = 490.70587
$ws.Range("L80").Value = 490.70587
$ws.Range("N80").Value = -2486.70587

$ws.Range("H83").Value = 412.0435
$ws.Range("J83").Value = 490.70587
$ws.Range("L83").Value = 2453.52935
$ws.Range("N83").Value = -12437.52935

$ws.Range("H105").Value = 2236.7693
$ws.Range("I105").Value = 1372.6786
$ws.Range("J105").Value = 4436.273
$ws.Range("K105").Value = 1372.6786
$ws.Range("L105").Value = 4436.273
$ws.Range("M105").Value = 374.3214
$ws.Range("N105").Value = -7930.273

$ws.Range("H107").Value = 4710.2085
$ws.Range("I107").Value = 5418.9414
$ws.Range("J107").Value = 2989
$ws.Range("K107").Value = 5418.9414
$ws.Range("L107").Value = 2989
$ws.Range("M107").Value = -3498.9414
$ws.Range("N107").Value = -6829

$ws.Range("H134").Value = 61406.76
$ws.Range("I134").Value = 82315.57000000001
$ws.Range("J134").Value = 1897.0769
$ws.Range("K134").Value = 246946.71
$ws.Range("L134").Value = 5691.2307
$ws.Range("M134").Value = -244411.71
$ws.Range("N134").Value = -10761.2307

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3873.5
$ws.Range("I16").Value = 4208.2
$ws.Range("J16").Value = 2200
$ws.Range("K16").Value = 4208.2
$ws.Range("L16").Value = 2200
$ws.Range("M16").Value = -3921.2
$ws.Range("N16").Value = -2774

$ws.Range("H31").Value = 1897.7727
$ws.Range("I31").Value = 1442.2258
$ws.Range("K31").Value = 1442.2258
$ws.Range("M31").Value = -1147.2258

$ws.Range("H34").Value = 1897.7727
$ws.Range("I34").Value = 1442.2258
$ws.Range("K34").Value = 1442.2258
$ws.Range("M34").Value = -1240.2258

$ws.Range("H113").Value = 3873.5
$ws.Range("I113").Value = 4208.2
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 4208.2
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = -2038.2
$ws.Range("N113").Value = -6540

$ws.Range("H132").Value = 1908.7959
$ws.Range("I132").Value = 1468.25
$ws.Range("J132").Value = 3866.7778
$ws.Range("K132").Value = 4404.75
$ws.Range("L132").Value = 11600.3334
$ws.Range("M132").Value = -1874.75
$ws.Range("N132").Value = -16660.3334

$ws.Range("H134").Value = 4915.278
$ws.Range("I134").Value = 5094.143
$ws.Range("K134").Value = 15282.429
$ws.Range("M134").Value = -12747.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 912
$ws.Range("I131").Value = 492
$ws.Range("J131").Value = 934.8261
$ws.Range("K131").Value = 1476
$ws.Range("L131").Value = 2804.4783
$ws.Range("M131").Value = 3564
$ws.Range("N131").Value = -12884.4783

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 27333.152
$ws.Range("I70").Value = 28336.363
$ws.Range("J70").Value = 5262.5
$ws.Range("K70").Value = 28336.363
$ws.Range("L70").Value = 5262.5
$ws.Range("M70").Value = -28066.363
$ws.Range("N70").Value = -5802.5

$ws.Range("H73").Value = 27333.152
$ws.Range("I73").Value = 28336.363
$ws.Range("J73").Value = 5262.5
$ws.Range("K73").Value = 28336.363
$ws.Range("L73").Value = 5262.5
$ws.Range("M73").Value = -27400.363
$ws.Range("N73").Value = -7134.5

$ws.Range("H122").Value = 2059.7273
$ws.Range("I122").Value = 2287.7144
$ws.Range("J122").Value = 1953.3334
$ws.Range("K122").Value = 6863.1432
$ws.Range("L122").Value = 5860.0002
$ws.Range("M122").Value = -4413.1432
$ws.Range("N122").Value = -10760.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 710.7778
$ws.Range("I96").Value = 671.8570999999999
$ws.Range("J96").Value = 847
$ws.Range("K96").Value = 671.8570999999999
$ws.Range("L96").Value = 847
$ws.Range("M96").Value = 701.1429000000001
$ws.Range("N96").Value = -3593

$ws.Range("H122").Value = 8837.759
$ws.Range("I122").Value = 10252
$ws.Range("J122").Value = 3416.5
$ws.Range("K122").Value = 30756
$ws.Range("L122").Value = 10249.5
$ws.Range("M122").Value = -28306
$ws.Range("N122").Value = -15149.5
